$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2222.5557  # H86: was 2230.4
$ws.Cells.Item(86, 9).Value = 2057.7144  # I86: was 2215
$ws.Cells.Item(86, 10).Value = 2799.5  # J86: was 2266.3333
$ws.Cells.Item(86, 11).Value = 2057.7144  # K86: was 2215
$ws.Cells.Item(86, 12).Value = 2799.5  # L86: was 2266.3333
$ws.Cells.Item(86, 13).Value = -934.7143999999998  # M86: was -1092
$ws.Cells.Item(86, 14).Value = -5045.5  # N86: was -4512.3333

$ws.Cells.Item(88, 8).Value = 0  # H88: was 4220
$ws.Cells.Item(88, 9).Value = 0  # I88: was 4990
$ws.Cells.Item(88, 10).Value = 0  # J88: was 3450
$ws.Cells.Item(88, 11).Value = 0  # K88: was 4990
$ws.Cells.Item(88, 12).Value = 0  # L88: was 3450
$ws.Cells.Item(88, 13).ClearContents()  # M88: was -4584
$ws.Cells.Item(88, 14).ClearContents()  # N88: was -4262

$ws.Cells.Item(89, 8).Value = 2222.5557  # H89: was 2230.4
$ws.Cells.Item(89, 9).Value = 2057.7144  # I89: was 2215
$ws.Cells.Item(89, 10).Value = 2799.5  # J89: was 2266.3333
$ws.Cells.Item(89, 11).Value = 10288.572  # K89: was 11075
$ws.Cells.Item(89, 12).Value = 13997.5  # L89: was 11331.6665
$ws.Cells.Item(89, 13).Value = -4672.572  # M89: was -5459
$ws.Cells.Item(89, 14).Value = -25229.5  # N89: was -22563.6665

$ws.Cells.Item(91, 8).Value = 0  # H91: was 4220
$ws.Cells.Item(91, 9).Value = 0  # I91: was 4990
$ws.Cells.Item(91, 10).Value = 0  # J91: was 3450
$ws.Cells.Item(91, 11).Value = 0  # K91: was 4990
$ws.Cells.Item(91, 12).Value = 0  # L91: was 3450
$ws.Cells.Item(91, 13).ClearContents()  # M91: was -3586
$ws.Cells.Item(91, 14).ClearContents()  # N91: was -6258

$ws.Cells.Item(116, 8).Value = 11552.5  # H116: was 9284.833000000001
$ws.Cells.Item(116, 9).Value = 5999.5  # I116: was 5374.5
$ws.Cells.Item(116, 11).Value = 5999.5  # K116: was 5374.5
$ws.Cells.Item(116, 13).Value = -2557.5  # M116: was -1932.5

$ws.Cells.Item(125, 8).Value = 19999  # H125: was 10499.5
$ws.Cells.Item(125, 10).Value = 19999  # J125: was 10499.5
$ws.Cells.Item(125, 12).Value = 179991  # L125: was 94495.5
$ws.Cells.Item(125, 14).Value = -184911  # N125: was -99415.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1399  # H61: was 1463.6364
$ws.Cells.Item(61, 9).Value = 1399  # I61: was 1463.6364
$ws.Cells.Item(61, 11).Value = 1399  # K61: was 1463.6364
$ws.Cells.Item(61, 13).Value = -1187  # M61: was -1251.6364

$ws.Cells.Item(63, 8).Value = 3797.55  # H63: was 3792.8
$ws.Cells.Item(63, 9).Value = 2438.3235  # I63: was 2432.7354
$ws.Cells.Item(63, 11).Value = 2438.3235  # K63: was 2432.7354
$ws.Cells.Item(63, 13).Value = -1752.3235  # M63: was -1746.7354

$ws.Cells.Item(66, 8).Value = 3797.55  # H66: was 3792.8
$ws.Cells.Item(66, 9).Value = 2438.3235  # I66: was 2432.7354
$ws.Cells.Item(66, 11).Value = 12191.6175  # K66: was 12163.677
$ws.Cells.Item(66, 13).Value = -8759.6175  # M66: was -8731.677

$ws.Cells.Item(74, 8).Value = 1663.3529  # H74: was 1775.3125
$ws.Cells.Item(74, 9).Value = 1418.6  # I74: was 1519.7693
$ws.Cells.Item(74, 10).Value = 3499  # J74: was 2882.6667
$ws.Cells.Item(74, 11).Value = 1418.6  # K74: was 1519.7693
$ws.Cells.Item(74, 12).Value = 3499  # L74: was 2882.6667
$ws.Cells.Item(74, 13).Value = -544.5999999999999  # M74: was -645.7692999999999
$ws.Cells.Item(74, 14).Value = -5247  # N74: was -4630.6667

$ws.Cells.Item(77, 8).Value = 1663.3529  # H77: was 1775.3125
$ws.Cells.Item(77, 9).Value = 1418.6  # I77: was 1519.7693
$ws.Cells.Item(77, 10).Value = 3499  # J77: was 2882.6667
$ws.Cells.Item(77, 11).Value = 7093  # K77: was 7598.8465
$ws.Cells.Item(77, 12).Value = 17495  # L77: was 14413.3335
$ws.Cells.Item(77, 13).Value = -2725  # M77: was -3230.8465
$ws.Cells.Item(77, 14).Value = -26231  # N77: was -23149.3335

$ws.Cells.Item(95, 8).Value = 53600.855  # H95: was 55743.715
$ws.Cells.Item(95, 10).Value = 53600.855  # J95: was 55743.715
$ws.Cells.Item(95, 12).Value = 53600.855  # L95: was 55743.715
$ws.Cells.Item(95, 14).Value = -59092.855  # N95: was -61235.715

$ws.Cells.Item(136, 8).Value = 1399  # H136: was 1463.6364
$ws.Cells.Item(136, 9).Value = 1399  # I136: was 1463.6364
$ws.Cells.Item(136, 11).Value = 4197  # K136: was 4390.9092
$ws.Cells.Item(136, 13).Value = -1647  # M136: was -1840.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(17, 8).Value = 6500  # H17: was 3254.5
$ws.Cells.Item(17, 10).Value = 6500  # J17: was 3254.5
$ws.Cells.Item(17, 12).Value = 6500  # L17: was 3254.5
$ws.Cells.Item(17, 14).Value = -6844  # N17: was -3598.5

$ws.Cells.Item(86, 8).Value = 3430.3333  # H86: was 3441.389
$ws.Cells.Item(86, 9).Value = 1265.375  # I86: was 1277.8125
$ws.Cells.Item(86, 11).Value = 1265.375  # K86: was 1277.8125
$ws.Cells.Item(86, 13).Value = -142.375  # M86: was -154.8125

$ws.Cells.Item(89, 8).Value = 3430.3333  # H89: was 3441.389
$ws.Cells.Item(89, 9).Value = 1265.375  # I89: was 1277.8125
$ws.Cells.Item(89, 11).Value = 6326.875  # K89: was 6389.0625
$ws.Cells.Item(89, 13).Value = -710.875  # M89: was -773.0625

$ws.Cells.Item(107, 8).Value = 3118.1667  # H107: was 1811.6666
$ws.Cells.Item(107, 9).Value = 3118.1667  # I107: was 1966.2727
$ws.Cells.Item(107, 10).Value = 0  # J107: was 111
$ws.Cells.Item(107, 11).Value = 3118.1667  # K107: was 1966.2727
$ws.Cells.Item(107, 12).Value = 0  # L107: was 111
$ws.Cells.Item(107, 13).Value = -1198.1667  # M107: was -46.27269999999999
$ws.Cells.Item(107, 14).ClearContents()  # N107: was -3951

$ws.Cells.Item(134, 8).Value = 1052.5  # H134: was 1299
$ws.Cells.Item(134, 9).Value = 1052.5  # I134: was 1299
$ws.Cells.Item(134, 11).Value = 3157.5  # K134: was 3897
$ws.Cells.Item(134, 13).Value = -622.5  # M134: was -1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 0  # H31: was 4990
$ws.Cells.Item(31, 9).Value = 0  # I31: was 4990
$ws.Cells.Item(31, 11).Value = 0  # K31: was 4990
$ws.Cells.Item(31, 13).ClearContents()  # M31: was -4695

$ws.Cells.Item(34, 8).Value = 0  # H34: was 4990
$ws.Cells.Item(34, 9).Value = 0  # I34: was 4990
$ws.Cells.Item(34, 11).Value = 0  # K34: was 4990
$ws.Cells.Item(34, 13).ClearContents()  # M34: was -4788

$ws.Cells.Item(132, 8).Value = 1683.25  # H132: was 1695.3334
$ws.Cells.Item(132, 9).Value = 1683.25  # I132: was 1695.3334
$ws.Cells.Item(132, 11).Value = 5049.75  # K132: was 5086.0002
$ws.Cells.Item(132, 13).Value = -2519.75  # M132: was -2556.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 2120.3  # H60: was 2564
$ws.Cells.Item(60, 9).Value = 623.2  # I60: was 862.6667
$ws.Cells.Item(60, 10).Value = 3617.4  # J60: was 3840
$ws.Cells.Item(60, 11).Value = 1869.6  # K60: was 2588.0001
$ws.Cells.Item(60, 12).Value = 10852.2  # L60: was 11520
$ws.Cells.Item(60, 13).Value = -1618.6  # M60: was -2337.0001
$ws.Cells.Item(60, 14).Value = -11354.2  # N60: was -12022

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 256.77777  # H2: was 214.36363
$ws.Cells.Item(2, 9).Value = 247.25  # I2: was 202.5
$ws.Cells.Item(2, 11).Value = 247.25  # K2: was 202.5
$ws.Cells.Item(2, 13).Value = -134.25  # M2: was -89.5

$ws.Cells.Item(12, 8).Value = 806312.4  # H12: was 645050.4
$ws.Cells.Item(12, 9).Value = 1286666.4  # I12: was 772000.8
$ws.Cells.Item(12, 11).Value = 1286666.4  # K12: was 772000.8
$ws.Cells.Item(12, 13).Value = -1286526.4  # M12: was -771860.8

$ws.Cells.Item(25, 8).Value = 0  # H25: was 809
$ws.Cells.Item(25, 10).Value = 0  # J25: was 809
$ws.Cells.Item(25, 12).Value = 0  # L25: was 809
$ws.Cells.Item(25, 14).ClearContents()  # N25: was -1867

$ws.Cells.Item(109, 8).Value = 0  # H109: was 44996
$ws.Cells.Item(109, 10).Value = 0  # J109: was 44996
$ws.Cells.Item(109, 12).Value = 0  # L109: was 44996
$ws.Cells.Item(109, 14).ClearContents()  # N109: was -47076

$ws.Cells.Item(122, 8).Value = 2557.4707  # H122: was 2811.0667
$ws.Cells.Item(122, 9).Value = 1822.5834  # I122: was 2024.3636
$ws.Cells.Item(122, 10).Value = 4321.2  # J122: was 4974.5
$ws.Cells.Item(122, 11).Value = 5467.7502  # K122: was 6073.0908
$ws.Cells.Item(122, 12).Value = 12963.6  # L122: was 14923.5
$ws.Cells.Item(122, 13).Value = -3017.7502  # M122: was -3623.0908
$ws.Cells.Item(122, 14).Value = -17863.6  # N122: was -19823.5

$ws.Cells.Item(126, 8).Value = 2810.6667  # H126: was 2732.923
$ws.Cells.Item(126, 9).Value = 2810.6667  # I126: was 2732.923
$ws.Cells.Item(126, 11).Value = 8432.000100000001  # K126: was 8198.769
$ws.Cells.Item(126, 13).Value = -5962.000100000001  # M126: was -5728.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1902.6875  # H22: was 2294.9
$ws.Cells.Item(22, 9).Value = 1742.8572  # I22: was 2300
$ws.Cells.Item(22, 10).Value = 2027  # J22: was 2291.5
$ws.Cells.Item(22, 11).Value = 1742.8572  # K22: was 2300
$ws.Cells.Item(22, 12).Value = 2027  # L22: was 2291.5
$ws.Cells.Item(22, 13).Value = -1447.8572  # M22: was -2005
$ws.Cells.Item(22, 14).Value = -2617  # N22: was -2881.5

$ws.Cells.Item(27, 8).Value = 1902.6875  # H27: was 2294.9
$ws.Cells.Item(27, 9).Value = 1742.8572  # I27: was 2300
$ws.Cells.Item(27, 10).Value = 2027  # J27: was 2291.5
$ws.Cells.Item(27, 11).Value = 1742.8572  # K27: was 2300
$ws.Cells.Item(27, 12).Value = 2027  # L27: was 2291.5
$ws.Cells.Item(27, 13).Value = -1635.8572  # M27: was -2193
$ws.Cells.Item(27, 14).Value = -2241  # N27: was -2505.5

$ws.Cells.Item(46, 8).Value = 2999.8572  # H46: was 3100
$ws.Cells.Item(46, 9).Value = 2799.8  # I46: was 2800
$ws.Cells.Item(46, 11).Value = 2799.8  # K46: was 2800
$ws.Cells.Item(46, 13).Value = -2611.8  # M46: was -2612

$ws.Cells.Item(55, 8).Value = 887.8946999999999  # H55: was 982.94116
$ws.Cells.Item(55, 9).Value = 559.8889  # I55: was 621.25
$ws.Cells.Item(55, 10).Value = 1183.1  # J55: was 1304.4445
$ws.Cells.Item(55, 11).Value = 559.8889  # K55: was 621.25
$ws.Cells.Item(55, 12).Value = 1183.1  # L55: was 1304.4445
$ws.Cells.Item(55, 13).Value = -386.8889  # M55: was -448.25
$ws.Cells.Item(55, 14).Value = -1529.1  # N55: was -1650.4445

$ws.Cells.Item(93, 8).Value = 1234  # H93: was 1317
$ws.Cells.Item(93, 9).Value = 1234  # I93: was 1317
$ws.Cells.Item(93, 11).Value = 1234  # K93: was 1317
$ws.Cells.Item(93, 13).Value = 14  # M93: was -69

$ws.Cells.Item(100, 8).Value = 3379.8  # H100: was 3974.75
$ws.Cells.Item(100, 9).Value = 2474.75  # I100: was 2966.3333
$ws.Cells.Item(100, 11).Value = 2474.75  # K100: was 2966.3333
$ws.Cells.Item(100, 13).Value = -1933.75  # M100: was -2425.3333

$ws.Cells.Item(130, 8).Value = 84950  # H130: was 0
$ws.Cells.Item(130, 10).Value = 84950  # J130: was 0
$ws.Cells.Item(130, 12).Value = 84950  # L130: was 0
$ws.Cells.Item(130, 14).Value = -94990  # N130: was None

$ws.Cells.Item(132, 8).Value = 2374.75  # H132: was 2449.3157
$ws.Cells.Item(132, 9).Value = 2294.5789  # I132: was 2393.7058
$ws.Cells.Item(132, 10).Value = 3898  # J132: was 2922
$ws.Cells.Item(132, 11).Value = 6883.736699999999  # K132: was 7181.117400000001
$ws.Cells.Item(132, 12).Value = 11694  # L132: was 8766
$ws.Cells.Item(132, 13).Value = -4353.736699999999  # M132: was -4651.117400000001
$ws.Cells.Item(132, 14).Value = -16754  # N132: was -13826

$ws.Cells.Item(133, 8).Value = 0  # H133: was 90000
$ws.Cells.Item(133, 10).Value = 0  # J133: was 90000
$ws.Cells.Item(133, 12).Value = 0  # L133: was 90000
$ws.Cells.Item(133, 14).ClearContents()  # N133: was -95060

$ws.Cells.Item(136, 8).Value = 3944.3572  # H136: was 4132.5386
$ws.Cells.Item(136, 9).Value = 4812.7  # I136: was 4832.6
$ws.Cells.Item(136, 10).Value = 1773.5  # J136: was 1799
$ws.Cells.Item(136, 11).Value = 14438.1  # K136: was 14497.8
$ws.Cells.Item(136, 12).Value = 5320.5  # L136: was 5397
$ws.Cells.Item(136, 13).Value = -11888.1  # M136: was -11947.8
$ws.Cells.Item(136, 14).Value = -10420.5  # N136: was -10497

$ws.Cells.Item(139, 8).Value = 15000  # H139: was 0
$ws.Cells.Item(139, 9).Value = 15000  # I139: was 0
$ws.Cells.Item(139, 11).Value = 15000  # K139: was 0
$ws.Cells.Item(139, 13).Value = -9860  # M139: was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4691.4614  # H81: was 4791.4614
$ws.Cells.Item(81, 9).Value = 3742.5715  # I81: was 4133
$ws.Cells.Item(81, 10).Value = 5798.5  # J81: was 5355.857
$ws.Cells.Item(81, 11).Value = 7485.143  # K81: was 8266
$ws.Cells.Item(81, 12).Value = 11597  # L81: was 10711.714
$ws.Cells.Item(81, 13).Value = -6424.143  # M81: was -7205
$ws.Cells.Item(81, 14).Value = -13719  # N81: was -12833.714

$ws.Cells.Item(84, 8).Value = 4691.4614  # H84: was 4791.4614
$ws.Cells.Item(84, 9).Value = 3742.5715  # I84: was 4133
$ws.Cells.Item(84, 10).Value = 5798.5  # J84: was 5355.857
$ws.Cells.Item(84, 11).Value = 37425.715  # K84: was 41330
$ws.Cells.Item(84, 12).Value = 57985  # L84: was 53558.57
$ws.Cells.Item(84, 13).Value = -32121.715  # M84: was -36026
$ws.Cells.Item(84, 14).Value = -68593  # N84: was -64166.57

$ws.Cells.Item(132, 8).Value = 2734.3635  # H132: was 2337.9285
$ws.Cells.Item(132, 9).Value = 1109.7142  # I132: was 1042.1
$ws.Cells.Item(132, 11).Value = 3329.1426  # K132: was 3126.3
$ws.Cells.Item(132, 13).Value = -799.1425999999997  # M132: was -596.2999999999997

$ws.Cells.Item(136, 8).Value = 3884.3447  # H136: was 3484.6667
$ws.Cells.Item(136, 9).Value = 4013.1538  # I136: was 3374.625
$ws.Cells.Item(136, 10).Value = 3779.6875  # J136: was 3588.2354
$ws.Cells.Item(136, 11).Value = 12039.4614  # K136: was 10123.875
$ws.Cells.Item(136, 12).Value = 11339.0625  # L136: was 10764.7062
$ws.Cells.Item(136, 13).Value = -9489.4614  # M136: was -7573.875
$ws.Cells.Item(136, 14).Value = -16439.0625  # N136: was -15864.7062
